# Weekly update: a new price record for "Ciboulette" at Vega Modelo de Temuco
# is prepended to the data block that starts at row 148 (the first data row
# for this variety). All rows from the old 148 down to 275 shift down by one
# row (148->149 ... 275->276); the new row 148 carries the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 148, pushing everything below it down one row.
$ws.Rows.Item(148).EntireRow.Insert()

# Populate the newly inserted row 148 with this week's record.
$ws.Cells.Item(148, 1).Value  = 10                                # Mercado ID
$ws.Cells.Item(148, 2).Value  = "Vega Modelo de Temuco"           # Mercado
$ws.Cells.Item(148, 3).Value  = "La Araucanía"                    # Región
$ws.Cells.Item(148, 4).Value  = 44790                             # Fecha
$ws.Cells.Item(148, 5).Value  = 9                                 # Codreg
$ws.Cells.Item(148, 6).Value  = 100112039                         # Categoría ID
$ws.Cells.Item(148, 7).Value  = "Ciboulette"                      # Categoría
$ws.Cells.Item(148, 8).Value  = "Sin especificar"                 # Variedad
$ws.Cells.Item(148, 9).Value  = "Primera"                         # Calidad
$ws.Cells.Item(148, 10).Value = 55                                # Volumen
$ws.Cells.Item(148, 11).Value = 7000                              # Precio mínimo
$ws.Cells.Item(148, 12).Value = 7000                              # Precio máximo
$ws.Cells.Item(148, 13).Value = 7000                              # Precio promedio ponderado
$ws.Cells.Item(148, 14).Value = "$/docena de atados"              # Unidad de comercialización
$ws.Cells.Item(148, 15).Value = "Provincia de Cautín"             # Origen
$ws.Cells.Item(148, 16).Value = 2333                              # Precio $/Kg
$ws.Cells.Item(148, 17).Value = 3                                 # Kg o Unidades
$ws.Cells.Item(148, 18).Value = "Hortaliza"                       # Clasificación
